{"js": "// Update Leave Card \u2014 the cached MERGEFIELD result runs (employee\n// salutation/name/position/office) and the computed leave amounts are\n// replaced with the new employee's data.\nconst body = context.document.body;\n\nasync function replaceAll(findText, replaceText) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replaceText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// Salutation: MR -> MS\nawait replaceAll(\"MR\", \"MS\");\n\n// Employee full name\nawait replaceAll(\"MARIO A. VILLANUEVA\", \"MA. CARMELA D. ARELLANO\");\n\n// Position\nawait replaceAll(\"Admin Aide III\", \"Casual Medical Technologist\");\n\n// Office\nawait replaceAll(\"TOPS Detailed At Hanggang Sa Kabilang Buhay Services\", \"Ospital Ng Tagaytay\");\n\n// Monthly salary (appears twice: \"Monthly Salary\" line and \"TLB = S x D x CF\" line)\nawait replaceAll(\"13,522.00\", \"24,167.00\");\n\n// Total leave credits (appears twice: \"Total Leave Credits\" line and \"TLB = S x D x CF\" line)\nawait replaceAll(\"287.375\", \"102.500\");\n\n// Total leave benefits\nawait replaceAll(\"187,271.28\", \"119,378.98\");\n", "ps1": "# Update Leave Card \u2014 replace mail-merge result text (cached MERGEFIELD\n# values) with the new employee's data, and the recomputed leave amounts.\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute(\n        $findText,       # FindText\n        $true,           # MatchCase\n        $true,           # MatchWholeWord\n        $false,          # MatchWildcards\n        $false,          # MatchSoundsLike\n        $false,          # MatchAllWordForms\n        $true,           # Forward\n        1,               # Wrap (wdFindContinue)\n        $false,          # Format\n        $replaceText,    # ReplaceWith\n        2                # Replace (wdReplaceAll)\n    ) | Out-Null\n}\n\n# Salutation: MR -> MS\nReplace-Text \"MR\" \"MS\"\n\n# Employee name\nReplace-Text \"MARIO A. VILLANUEVA\" \"MA. CARMELA D. ARELLANO\"\n\n# Position\nReplace-Text \"Admin Aide III\" \"Casual Medical Technologist\"\n\n# Office\nReplace-Text \"TOPS Detailed At Hanggang Sa Kabilang Buhay Services\" \"Ospital Ng Tagaytay\"\n\n# Monthly salary (appears twice: \"Monthly Salary\" line and \"TLB = S x D x CF\" line)\nReplace-Text \"13,522.00\" \"24,167.00\"\n\n# Total leave credits (appears twice: \"Total Leave Credits\" line and \"TLB = S x D x CF\" line)\nReplace-Text \"287.375\" \"102.500\"\n\n# Total leave benefits\nReplace-Text \"187,271.28\" \"119,378.98\"\n"}
